$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 99: Rumor Has It / Commanding Craftsman's Tea
$ws.Cells.Item(99, 8).Value = 535.6667
$ws.Cells.Item(99, 9).Value = 495.27274
$ws.Cells.Item(99, 10).Value = 980
$ws.Cells.Item(99, 11).Value = 1485.81822
$ws.Cells.Item(99, 12).Value = 2940
$ws.Cells.Item(99, 13).Value = 12.18177999999989
$ws.Cells.Item(99, 14).Value = -5936

# Row 107: Another Man's Ink / Enchanted Truegold Ink
$ws.Cells.Item(107, 8).Value = 9259722
$ws.Cells.Item(107, 10).Value = 1450
$ws.Cells.Item(107, 12).Value = 1450
$ws.Cells.Item(107, 14).Value = -5290

# Row 125: Body over Mind / Grade 5 Dexterity Alkahest
$ws.Cells.Item(125, 8).Value = 4074.389
$ws.Cells.Item(125, 9).Value = 12766
$ws.Cells.Item(125, 10).Value = 2987.9375
$ws.Cells.Item(125, 11).Value = 114894
$ws.Cells.Item(125, 12).Value = 26891.4375
$ws.Cells.Item(125, 13).Value = -112434
$ws.Cells.Item(125, 14).Value = -31811.4375

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Cells.Item(137, 8).Value = 1341.1818
$ws.Cells.Item(137, 9).Value = 1255.258
$ws.Cells.Item(137, 10).Value = 1546.0769
$ws.Cells.Item(137, 11).Value = 3765.774
$ws.Cells.Item(137, 12).Value = 4638.2307
$ws.Cells.Item(137, 13).Value = -1215.774
$ws.Cells.Item(137, 14).Value = -9738.2307

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Cells.Item(141, 8).Value = 1963.9166
$ws.Cells.Item(141, 9).Value = 1687.7858
$ws.Cells.Item(141, 10).Value = 2930.375
$ws.Cells.Item(141, 11).Value = 5063.357400000001
$ws.Cells.Item(141, 12).Value = 8791.125
$ws.Cells.Item(141, 13).Value = 116.6425999999992
$ws.Cells.Item(141, 14).Value = -19151.125

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Cells.Item(2, 8).Value = 1561.8182
$ws.Cells.Item(2, 9).Value = 1020
$ws.Cells.Item(2, 11).Value = 1020
$ws.Cells.Item(2, 13).Value = -907

# Row 32: Ingot We Trust / Steel Ingot
$ws.Cells.Item(32, 8).Value = 7753.852
$ws.Cells.Item(32, 9).Value = 5815.7334
$ws.Cells.Item(32, 11).Value = 5815.7334
$ws.Cells.Item(32, 13).Value = -5528.7334

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Cells.Item(61, 8).Value = 2791.45
$ws.Cells.Item(61, 9).Value = 2796.3276
$ws.Cells.Item(61, 10).Value = 2650
$ws.Cells.Item(61, 11).Value = 2796.3276
$ws.Cells.Item(61, 12).Value = 2650
$ws.Cells.Item(61, 13).Value = -2584.3276
$ws.Cells.Item(61, 14).Value = -3074

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Cells.Item(74, 8).Value = 1123.5333
$ws.Cells.Item(74, 9).Value = 977.6389
$ws.Cells.Item(74, 10).Value = 1707.1111
$ws.Cells.Item(74, 11).Value = 977.6389
$ws.Cells.Item(74, 12).Value = 1707.1111
$ws.Cells.Item(74, 13).Value = -103.6389
$ws.Cells.Item(74, 14).Value = -3455.1111

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Cells.Item(77, 8).Value = 1123.5333
$ws.Cells.Item(77, 9).Value = 977.6389
$ws.Cells.Item(77, 10).Value = 1707.1111
$ws.Cells.Item(77, 11).Value = 4888.194500000001
$ws.Cells.Item(77, 12).Value = 8535.5555
$ws.Cells.Item(77, 13).Value = -520.1945000000005
$ws.Cells.Item(77, 14).Value = -17271.5555

# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Cells.Item(110, 8).Value = 633.0769
$ws.Cells.Item(110, 9).Value = 633.0769
$ws.Cells.Item(110, 11).Value = 633.0769
$ws.Cells.Item(110, 13).Value = 1411.9231

# Row 116: No Scope / Titanbronze Ingot
$ws.Cells.Item(116, 8).Value = 1561.8182
$ws.Cells.Item(116, 9).Value = 1020
$ws.Cells.Item(116, 11).Value = 1020
$ws.Cells.Item(116, 13).Value = 1274

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Cells.Item(132, 8).Value = 1904.409
$ws.Cells.Item(132, 9).Value = 1076.6346
$ws.Cells.Item(132, 10).Value = 4979
$ws.Cells.Item(132, 11).Value = 3229.9038
$ws.Cells.Item(132, 12).Value = 14937
$ws.Cells.Item(132, 13).Value = -699.9038
$ws.Cells.Item(132, 14).Value = -19997

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Cells.Item(136, 8).Value = 2791.45
$ws.Cells.Item(136, 9).Value = 2796.3276
$ws.Cells.Item(136, 10).Value = 2650
$ws.Cells.Item(136, 11).Value = 8388.9828
$ws.Cells.Item(136, 12).Value = 7950
$ws.Cells.Item(136, 13).Value = -5838.9828
$ws.Cells.Item(136, 14).Value = -13050

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Cells.Item(3, 8).Value = 1561.8182
$ws.Cells.Item(3, 9).Value = 1020
$ws.Cells.Item(3, 11).Value = 1020
$ws.Cells.Item(3, 13).Value = -906

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Cells.Item(107, 8).Value = 2006.5
$ws.Cells.Item(107, 9).Value = 3000
$ws.Cells.Item(107, 10).Value = 1013
$ws.Cells.Item(107, 11).Value = 3000
$ws.Cells.Item(107, 12).Value = 1013
$ws.Cells.Item(107, 13).Value = -1080
$ws.Cells.Item(107, 14).Value = -4853

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Cells.Item(134, 8).Value = 3578.3635
$ws.Cells.Item(134, 9).Value = 3691.5217
$ws.Cells.Item(134, 11).Value = 11074.5651
$ws.Cells.Item(134, 13).Value = -8539.5651

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof / Ash Lumber
$ws.Cells.Item(16, 8).Value = 1414.125
$ws.Cells.Item(16, 9).Value = 1141.1666
$ws.Cells.Item(16, 10).Value = 1577.9
$ws.Cells.Item(16, 11).Value = 1141.1666
$ws.Cells.Item(16, 12).Value = 1577.9
$ws.Cells.Item(16, 13).Value = -854.1666
$ws.Cells.Item(16, 14).Value = -2151.9

# Row 31: Wall Not Found / Walnut Lumber
$ws.Cells.Item(31, 8).Value = 3105.02
$ws.Cells.Item(31, 9).Value = 1504.7354
$ws.Cells.Item(31, 10).Value = 6505.625
$ws.Cells.Item(31, 11).Value = 1504.7354
$ws.Cells.Item(31, 12).Value = 6505.625
$ws.Cells.Item(31, 13).Value = -1209.7354
$ws.Cells.Item(31, 14).Value = -7095.625

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Cells.Item(34, 8).Value = 3105.02
$ws.Cells.Item(34, 9).Value = 1504.7354
$ws.Cells.Item(34, 10).Value = 6505.625
$ws.Cells.Item(34, 11).Value = 1504.7354
$ws.Cells.Item(34, 12).Value = 6505.625
$ws.Cells.Item(34, 13).Value = -1302.7354
$ws.Cells.Item(34, 14).Value = -6909.625

# Row 107: Built to Last / White Oak Lumber
$ws.Cells.Item(107, 8).Value = 283.86273
$ws.Cells.Item(107, 9).Value = 346.13333
$ws.Cells.Item(107, 10).Value = 257.91666
$ws.Cells.Item(107, 11).Value = 346.13333
$ws.Cells.Item(107, 12).Value = 257.91666
$ws.Cells.Item(107, 13).Value = 1573.86667
$ws.Cells.Item(107, 14).Value = -4097.91666

# Row 113: Patient Patients / White Ash Lumber
$ws.Cells.Item(113, 8).Value = 1414.125
$ws.Cells.Item(113, 9).Value = 1141.1666
$ws.Cells.Item(113, 10).Value = 1577.9
$ws.Cells.Item(113, 11).Value = 1141.1666
$ws.Cells.Item(113, 12).Value = 1577.9
$ws.Cells.Item(113, 13).Value = 1028.8334
$ws.Cells.Item(113, 14).Value = -5917.9

# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Cells.Item(122, 8).Value = 2886
$ws.Cells.Item(122, 9).Value = 2103
$ws.Cells.Item(122, 10).Value = 4138.8
$ws.Cells.Item(122, 11).Value = 6309
$ws.Cells.Item(122, 12).Value = 12416.4
$ws.Cells.Item(122, 13).Value = -3859
$ws.Cells.Item(122, 14).Value = -17316.4

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Cells.Item(132, 8).Value = 1604.7
$ws.Cells.Item(132, 9).Value = 1171.9744
$ws.Cells.Item(132, 10).Value = 3138.9092
$ws.Cells.Item(132, 11).Value = 3515.9232
$ws.Cells.Item(132, 12).Value = 9416.7276
$ws.Cells.Item(132, 13).Value = -985.9232000000002
$ws.Cells.Item(132, 14).Value = -14476.7276

$ws = $wb.Worksheets.Item("CUL")
# Row 114: One Last Meal / Mushroom Saute
$ws.Cells.Item(114, 8).Value = 9205.076999999999
$ws.Cells.Item(114, 9).Value = 296.875
$ws.Cells.Item(114, 10).Value = 23458.2
$ws.Cells.Item(114, 11).Value = 890.625
$ws.Cells.Item(114, 12).Value = 70374.60000000001
$ws.Cells.Item(114, 13).Value = 2363.375
$ws.Cells.Item(114, 14).Value = -76882.60000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws.Cells.Item(107, 8).Value = 1460.8334
$ws.Cells.Item(107, 9).Value = 598.9167
$ws.Cells.Item(107, 10).Value = 2322.75
$ws.Cells.Item(107, 11).Value = 598.9167
$ws.Cells.Item(107, 12).Value = 2322.75
$ws.Cells.Item(107, 13).Value = 1321.0833
$ws.Cells.Item(107, 14).Value = -6162.75

# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Cells.Item(113, 8).Value = 55556856
$ws.Cells.Item(113, 9).Value = 90910160
$ws.Cells.Item(113, 10).Value = 1655.7142
$ws.Cells.Item(113, 11).Value = 90910160
$ws.Cells.Item(113, 12).Value = 1655.7142
$ws.Cells.Item(113, 13).Value = -90907990
$ws.Cells.Item(113, 14).Value = -5995.7142

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Cells.Item(126, 8).Value = 7557.5557
$ws.Cells.Item(126, 9).Value = 9618.076999999999
$ws.Cells.Item(126, 10).Value = 2200.2
$ws.Cells.Item(126, 11).Value = 28854.231
$ws.Cells.Item(126, 12).Value = 6600.599999999999
$ws.Cells.Item(126, 13).Value = -26384.231
$ws.Cells.Item(126, 14).Value = -11540.6

$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly / Raptor Leather
$ws.Cells.Item(61, 8).Value = 3799.9333
$ws.Cells.Item(61, 9).Value = 2999
$ws.Cells.Item(61, 10).Value = 4333.8887
$ws.Cells.Item(61, 11).Value = 2999
$ws.Cells.Item(61, 12).Value = 4333.8887
$ws.Cells.Item(61, 13).Value = -2797
$ws.Cells.Item(61, 14).Value = -4737.8887

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Cells.Item(113, 8).Value = 3799.9333
$ws.Cells.Item(113, 9).Value = 2999
$ws.Cells.Item(113, 10).Value = 4333.8887
$ws.Cells.Item(113, 11).Value = 2999
$ws.Cells.Item(113, 12).Value = 4333.8887
$ws.Cells.Item(113, 13).Value = -829
$ws.Cells.Item(113, 14).Value = -8673.8887

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Cells.Item(132, 8).Value = 9968403
$ws.Cells.Item(132, 9).Value = 14328326
$ws.Cells.Item(132, 10).Value = 2864.3572
$ws.Cells.Item(132, 11).Value = 42984978
$ws.Cells.Item(132, 12).Value = 8593.071599999999
$ws.Cells.Item(132, 13).Value = -42982448
$ws.Cells.Item(132, 14).Value = -13653.0716

$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax / Bright Linen Yarn
$ws.Cells.Item(107, 8).Value = 250000130
$ws.Cells.Item(107, 9).Value = 333333470
$ws.Cells.Item(107, 10).Value = 150
$ws.Cells.Item(107, 11).Value = 1000000410
$ws.Cells.Item(107, 12).Value = 450
$ws.Cells.Item(107, 13).Value = -999998490
$ws.Cells.Item(107, 14).Value = -4290

# Row 113: A Tender Table / Pixie Floss
$ws.Cells.Item(113, 8).Value = 1572.3636
$ws.Cells.Item(113, 9).Value = 1500.4286
$ws.Cells.Item(113, 10).Value = 1698.25
$ws.Cells.Item(113, 11).Value = 4501.2858
$ws.Cells.Item(113, 12).Value = 5094.75
$ws.Cells.Item(113, 13).Value = -2331.2858
$ws.Cells.Item(113, 14).Value = -9434.75

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Cells.Item(132, 8).Value = 1015.20337
$ws.Cells.Item(132, 9).Value = 728.5306399999999
$ws.Cells.Item(132, 10).Value = 2419.9
$ws.Cells.Item(132, 11).Value = 2185.59192
$ws.Cells.Item(132, 12).Value = 7259.700000000001
$ws.Cells.Item(132, 13).Value = 344.4080800000002
$ws.Cells.Item(132, 14).Value = -12319.7

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Cells.Item(136, 8).Value = 1093.0869
$ws.Cells.Item(136, 9).Value = 631.5714
$ws.Cells.Item(136, 10).Value = 1480.76
$ws.Cells.Item(136, 11).Value = 1894.7142
$ws.Cells.Item(136, 12).Value = 4442.28
$ws.Cells.Item(136, 13).Value = 655.2857999999999
$ws.Cells.Item(136, 14).Value = -9542.279999999999
